$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 4294.964340500327
$ws.Cells.Item(3, 2).Value = 4926.153965455343
$ws.Cells.Item(4, 2).Value = 5464.408501268852
$ws.Cells.Item(5, 2).Value = 6056.764868294075
$ws.Cells.Item(6, 2).Value = 6833.226149065153
$ws.Cells.Item(7, 2).Value = 7358.258706741193
$ws.Cells.Item(8, 2).Value = 7560.823388780289
$ws.Cells.Item(9, 2).Value = 8156.464403688821
$ws.Cells.Item(10, 2).Value = 8490.618780374585
$ws.Cells.Item(11, 2).Value = 9083.484782114929
$ws.Cells.Item(12, 2).Value = 9313.896030017604
$ws.Cells.Item(13, 2).Value = 9502.365027749456
$ws.Cells.Item(14, 2).Value = 9991.299151908721
$ws.Cells.Item(15, 2).Value = 10287.82191723134
$ws.Cells.Item(16, 2).Value = 10617.72211794718
$ws.Cells.Item(17, 2).Value = 10806.93903147523
$ws.Cells.Item(18, 2).Value = 11076.81108858635
$ws.Cells.Item(19, 2).Value = 11301.06367270383
$ws.Cells.Item(20, 2).Value = 11618.31939991226
$ws.Cells.Item(21, 2).Value = 11765.0961080947
$ws.Cells.Item(22, 2).Value = 11894.41129315888
$ws.Cells.Item(23, 2).Value = 11968.66750596234
$ws.Cells.Item(24, 2).Value = 12198.39333331799
$ws.Cells.Item(25, 2).Value = 12444.63308485151
$ws.Cells.Item(26, 2).Value = 12739.27831815544
$ws.Cells.Item(27, 2).Value = 12834.92636584731
$ws.Cells.Item(28, 2).Value = 13039.99433250989
$ws.Cells.Item(29, 2).Value = 13210.43681948788
$ws.Cells.Item(30, 2).Value = 13355.84276198929
$ws.Cells.Item(31, 2).Value = 13478.52368072836
$ws.Cells.Item(32, 2).Value = 13679.00729914817
$ws.Cells.Item(33, 2).Value = 13752.45818072856
$ws.Cells.Item(34, 2).Value = 13979.6642935727
$ws.Cells.Item(35, 2).Value = 14107.25641726138
$ws.Cells.Item(36, 2).Value = 13958.81434771256
$ws.Cells.Item(37, 2).Value = 14079.94621813633
$ws.Cells.Item(38, 2).Value = 14240.6376546265
$ws.Cells.Item(39, 2).Value = 14317.23417890966
$ws.Cells.Item(40, 2).Value = 14488.99035059528
$ws.Cells.Item(41, 2).Value = 14578.97998375641
$ws.Cells.Item(42, 2).Value = 14722.45546914963
$ws.Cells.Item(43, 2).Value = 14805.14108839904
$ws.Cells.Item(44, 2).Value = 14884.05182910904
$ws.Cells.Item(45, 2).Value = 15238.22996678246
$ws.Cells.Item(46, 2).Value = 15297.92464207829
$ws.Cells.Item(47, 2).Value = 15403.49262767791
$ws.Cells.Item(48, 2).Value = 15576.79312480013
$ws.Cells.Item(49, 2).Value = 15561.50801447334
$ws.Cells.Item(50, 2).Value = 15575.54402581014
$ws.Cells.Item(51, 2).Value = 15695.4091117376
$ws.Cells.Item(52, 2).Value = 15709.20306581252
$ws.Cells.Item(53, 2).Value = 15831.61629122414
$ws.Cells.Item(54, 2).Value = 15853.93741127316
$ws.Cells.Item(55, 2).Value = 16047.25741812299
$ws.Cells.Item(56, 2).Value = 15934.67028225771
$ws.Cells.Item(57, 2).Value = 16091.04934136898
$ws.Cells.Item(58, 2).Value = 16094.30334947824
$ws.Cells.Item(59, 2).Value = 15966.86512221618
$ws.Cells.Item(60, 2).Value = 16060.77889280614
$ws.Cells.Item(61, 2).Value = 16081.81664056357
$ws.Cells.Item(62, 2).Value = 16136.94882170101
